$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34 (weekly update adds a new price record),
# pushing the existing rows 34-80 down to 35-81.
$ws.Rows(34).Insert()

# Populate the newly inserted row 34 with the new weekly record.
$ws.Range("A34").Value = 3
$ws.Range("B34").Value = "Femacal de La Calera"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44803
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = 100112035
$ws.Range("G34").Value = "Bruselas (repollito)"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 85
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 15500
$ws.Range("M34").Value = 15265
$ws.Range("N34").Value = "`$/malla 15 kilos"
$ws.Range("O34").Value = "Provincia de Quillota"
$ws.Range("P34").Value = 1018
$ws.Range("Q34").Value = 15
$ws.Range("R34").Value = "Hortaliza"
